$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update quantity figures (column C) - formulas in column E recalc automatically
$ws.Range("C9").Value = 20
$ws.Range("C14").Value = 50
$ws.Range("C31").Value = 2000
$ws.Range("C32").Value = 500
$ws.Range("C43").Value = 231901

# Update the view: scroll position and active selection cell
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("I37").Select()
